# Add a new "F" column to the SPFE-JP sheet that mirrors the password/ID
# column (B) for each card row — used elsewhere (e.g. an FH passlist).
# Also clears the date stamp out of B1 (keeping its style) and mirrors an
# (empty, same-styled) cell into F1 so the header row lines up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPFE-JP")

# Mirror the header cell's style (but not its value) into the new column.
$ws.Range("B1").Copy($ws.Range("F1"))
$ws.Range("B1").Value = $null
$ws.Range("F1").Value = $null

# Copy the ID numbers from column B down into the new column F.
$ws.Range("B2:B26").Copy($ws.Range("F2:F26"))

# New column F should look/behave like column B (same width/bestFit-ish).
$ws.Columns("F").ColumnWidth = $ws.Columns("B").ColumnWidth

# Restore the selection/view the workbook was left in after the edit
# (user finished by selecting the B:D data block, landing on D26).
$ws.Range("B2:D26").Select()
